# The 9 "lci_methodN" worksheets each hold one LCIA method's results: column A
# (rows 2-15) carries the method name (as a repeated label), column B the
# process name, and C:F the four numeric result columns for that method.
#
# This edit rotates the per-method content (label + numeric results) across
# the 9 sheets while the sheets themselves (tabs, B-column process names)
# stay put: sheet k's new content = sheet mapping[k]'s old content, i.e. a
# cyclic shift by two positions (sheet1<-sheet3, sheet2<-sheet4, ...,
# sheet7<-sheet9, sheet8<-sheet1, sheet9<-sheet2).

$wb = $excel.ActiveWorkbook

$labelAddr = "A2:A15"
$dataAddr  = "C2:F15"

# Snapshot every sheet's current (before) label + data first -- writes must
# not clobber a block that a later sheet still needs to read.
$labels = @{}
$data = @{}
for ($i = 1; $i -le 9; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $labels[$i] = $ws.Range($labelAddr).Cells.Item(1,1).Value()
    $data[$i] = $ws.Range($dataAddr).Value()
}

# new sheet k <- old sheet mapping[k]
$mapping = @{1=3; 2=4; 3=5; 4=6; 5=7; 6=8; 7=9; 8=1; 9=2}

for ($i = 1; $i -le 9; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $src = $mapping[$i]
    $ws.Range($labelAddr).Value = $labels[$src]
    $ws.Range($dataAddr).Value = $data[$src]
}
